$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-40 contain Id values "men_boohoo_hoodie_<n>" in column A.
# Rename them to "men_boohoohoodie_<n>" (drop the underscore between
# "boohoo" and "hoodie").
for ($row = 2; $row -le 40; $row++) {
    $n = $row - 1
    $ws.Cells.Item($row, 1).Value = "men_boohoohoodie_$n"
}
